$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hedge with zcb")

$values = @(
    -27872.73217226321,
    -28505.19157816322,
    -28951.16391042173,
    -29322.46222077663,
    -29465.67670647293,
    -29270.15557835265,
    -29085.13773181424,
    -28919.83322262217,
    -28608.05240017437,
    -28150.93408682003,
    -27664.48673328008,
    -27253.25510337381,
    -26723.8192678688,
    -26245.77455861241,
    -25256.98877275228,
    -24305.70329757581,
    -23973.6899340907,
    -23287.3223833589,
    -22683.44833445781,
    -22051.95055150932,
    -21393.72180435376,
    -20714.1215496283,
    -20140.64076975684,
    -126522.0942230288,
    -15966.05072521717,
    -15773.91933315351,
    -15472.58800625017,
    -14941.39418334994,
    -14709.02088501978,
    -14253.40069303035,
    -13963.84161733983,
    -13485.23445993945,
    -13194.01049688353,
    -12842.43416880642,
    -12468.07131298199,
    -118230.1708920357,
    -10061.68160939902,
    -9770.351655640088,
    -9576.845055425692,
    -9399.691138934173,
    -9482.200229550292,
    -9179.840758572129,
    -8818.489126861503,
    -8740.290892814754,
    -8532.524324760823,
    -8435.810227489135,
    -8287.447494908836,
    -8084.70816677213,
    -7858.009488631932,
    -7698.608058461469,
    -7626.789900368242,
    -7421.600177721641,
    -7354.08262542528,
    -7208.394319977569,
    -7004.423362970527,
    -6796.853457025359,
    -6689.817791154689,
    -6576.90674109983,
    -6471.391111436718,
    -375344.5150015376,
    -2501.172781160858,
    -2472.357800316604,
    -2426.488532461176,
    -2385.089875361895,
    -2342.260552554153,
    -2273.622500472011,
    -2246.925960998009,
    -2180.457705455039,
    -2171.495254655041,
    -2121.71188378467,
    -2113.768735606039,
    -2082.956201836494,
    -2016.789928209417,
    -2016.834970768178,
    -1997.291766057096,
    -2020.020294147966,
    -1948.731695737526,
    -1952.037018573134,
    -1880.039870115615,
    -1875.081761426459,
    -1907.631957430888,
    -1853.281170119549,
    -1835.790387995673,
    -120536.6535508656,
    -1069.16908048592,
    -1037.378309504954,
    -1006.000617172101,
    -959.7256246177544,
    -961.165293637821,
    -927.8541838848761,
    -936.8959503854447,
    -895.9969601476033,
    -911.9707307866365,
    -924.4015243295261,
    -905.1183041758055,
    -97405.6188670859,
    -576.2131538399767,
    -587.9169370934774,
    -576.4997619579231,
    -580.4754121272406,
    -584.211346050738,
    -549.2474733542732,
    -549.8961023201276,
    -564.1062702315463,
    -552.3092859922812,
    -571.077061912014,
    -564.2036018871468,
    -577.220006708626,
    -570.724020188372,
    -553.3339442735194,
    -530.3967074245548,
    -537.0728868010428,
    -550.1331633514874,
    -545.2585385317417,
    -519.8537324430229,
    -528.5157178316325,
    -508.5178836755719,
    -498.959295779198,
    -498.1610634008524,
    -111069.0481202427
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
